$d = $word.ActiveDocument
$tbl = $d.Tables(1)

function Insert-WingdingsCheck($rowIndex, $pOpenTag, $pPrXml) {
    $cell = $tbl.Cell($rowIndex, 2)
    $r = $cell.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        $pOpenTag + $pPrXml +
        '<w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="de-CH"/></w:rPr><w:sym w:font="Wingdings" w:char="F0FC"/></w:r>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# Row 5: "SettingsFrame: Usereingaben überprüfen"
Insert-WingdingsCheck 5 `
    '<w:p w:rsidR="007D6018" w:rsidRPr="007E1874" w:rsidRDefault="007D6018" w:rsidP="007D6018">' `
    '<w:pPr><w:cnfStyle w:val="000000010000"/><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr>'

# Row 12: "Field wechseln"
Insert-WingdingsCheck 12 `
    '<w:p w:rsidR="007D6018" w:rsidRPr="007E1874" w:rsidRDefault="007D6018" w:rsidP="007D6018">' `
    '<w:pPr><w:cnfStyle w:val="000000100000"/><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="de-CH"/></w:rPr></w:pPr>'

# Row 19: "Menukonzept: Menus deaktivieren"
Insert-WingdingsCheck 19 `
    '<w:p w:rsidR="00EF245F" w:rsidRPr="007E1874" w:rsidRDefault="00EF245F" w:rsidP="007D6018">' `
    '<w:pPr><w:cnfStyle w:val="000000010000"/><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="de-CH"/></w:rPr></w:pPr>'

Write-Host "Inserted Wingdings checkmarks into rows 5, 12, 19"
